$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column F
$ws.Range("F1").Value = "EDAM_DEF"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Populate EDAM_DEF values for rows 2-36
$ws.Range("F2").Value = '[''3D coordinate and associated data for a macromolecular tertiary (3D) structure or part of a structure.'']'
$ws.Range("F3").Value = '[''Specification of one or more colors.'']'
$ws.Range("F4").Value = '[''Alphabet for a DNA sequence with possible ambiguity, unknown positions and non-sequence characters.'']'
$ws.Range("F5").Value = '[''The cell cycle including key genes and proteins.'']'
$ws.Range("F6").Value = '[''The analysis of levels and patterns of synthesis of gene products (proteins and functional RNA) including interpretation in functional terms of gene expression data.'']'
$ws.Range("F7").Value = '[''Fungi and molds, e.g. information on a specific fungal genome including molecular sequences, genes and annotation.'']'
$ws.Range("F8").Value = '["Laboratory technique to sequence the complete DNA sequence of an organism''s genome at a single time."]'
$ws.Range("F9").Value = '[''The secondary structure assignment (predicted or real) of a nucleic acid or protein.'']'
$ws.Range("F10").Value = '[''Alphabet for a protein sequence with possible ambiguity, unknown positions and non-sequence characters.'']'
$ws.Range("F11").Value = '[''An analytical chemistry technique that measures the mass-to-charge ratio and abundance of ions in the gas phase.'']'
$ws.Range("F12").Value = '[''Alphabet for an RNA sequence with possible ambiguity, unknown positions and non-sequence characters.'']'
$ws.Range("F13").Value = '[''One or more ribosomal RNA (rRNA) sequences.'']'
$ws.Range("F14").Value = '[''One or more transfer RNA (tRNA) sequences.'']'
$ws.Range("F15").Value = '[''The study of matter by studying the interference pattern from firing electrons at a sample, to analyse structures at resolutions higher than can be achieved using light.'']'
$ws.Range("F16").Value = '[''Km is the concentration (usually in Molar units) of substrate that leads to half-maximal velocity of an enzyme-catalysed reaction.'']'
$ws.Range("F17").Value = '[''3D coordinate and associated data for a multi-protein complex; two or more polypeptides chains in a stable, functional association with one another.'']'
$ws.Range("F18").Value = '[''The maximum initial velocity or rate of a reaction. It is the limiting velocity as substrate concentrations get very large.'']'
$ws.Range("F19").Value = '[''Alphabet for a nucleotide sequence with possible ambiguity, unknown positions and non-sequence characters.'']'
$ws.Range("F20").Value = '[''File format of a CT (Connectivity Table) file from the RNAstructure package.'']'
$ws.Range("F21").Value = '[''Align molecular sequence to structure in 3D space (threading).'']'
$ws.Range("F22").Value = '[''Model or simulate protein-protein binding using comparative modelling or other techniques.'']'
$ws.Range("F23").Value = '[''Alignment (superimposition) of molecular tertiary (3D) structures.'']'
$ws.Range("F24").Value = '[''Alignment of multiple molecular sequences.'']'
$ws.Range("F25").Value = '[''Align more than two molecular sequences.'']'
$ws.Range("F26").Value = '[''Align exactly two molecular sequences.'']'
$ws.Range("F27").Value = '[''Search a tertiary structure database, typically by sequence and/or structure comparison, or some other means, and retrieve structures and associated data.'']'
$ws.Range("F28").Value = '[''A statistical procedure that uses an orthogonal transformation to convert a set of observations of possibly correlated variables into a set of values of linearly uncorrelated variables called principal components.'']'
$ws.Range("F29").Value = '[''Generate, process or analyse a biological pathway.'']'
$ws.Range("F30").Value = '[''The processing and analysis of natural language, such as scientific literature in English, in order to extract data and information, or to enable human-computer interaction.'']'
$ws.Range("F31").Value = '[''Generate, process or analyse a biological network.'']'
$ws.Range("F32").Value = '[''Virtual screening is used in drug discovery to identify potential drug compounds.  It involves searching libraries of small molecules in order to identify those molecules which are most likely to bind to a drug target (typically a protein receptor or enzyme).'']'
$ws.Range("F33").Value = '[''Model the structure of a protein in complex with a small molecule or another macromolecule.'']'
$ws.Range("F34").Value = '[''Imaging in sections (sectioning), through the use of a wave-generating device (tomograph) that generates an image (a tomogram).'']'
$ws.Range("F35").Value = '[''Predict antigenic determinant sites (epitopes) in protein sequences.'']'
$ws.Range("F36").Value = '[''Toxins and the adverse effects of these chemical substances on living organisms.'']'

$excel.CutCopyMode = 0
